$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2506.6667
$ws.Range("I51").Value = 2000
$ws.Range("J51").Value = 2570
$ws.Range("K51").Value = 2000
$ws.Range("L51").Value = 2570
$ws.Range("M51").Value = -1516
$ws.Range("N51").Value = -3538
$ws.Range("H62").Value = 8770
$ws.Range("J62").Value = 17525
$ws.Range("L62").Value = 17525
$ws.Range("N62").Value = -18773
$ws.Range("H65").Value = 8770
$ws.Range("J65").Value = 17525
$ws.Range("L65").Value = 87625
$ws.Range("N65").Value = -93865
$ws.Range("H138").Value = 2299.7422
$ws.Range("J138").Value = 2314.5054
$ws.Range("L138").Value = 6943.5162
$ws.Range("N138").Value = -17223.5162
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 64084
$ws.Range("I10").Value = 400
$ws.Range("J10").Value = 80005
$ws.Range("K10").Value = 400
$ws.Range("L10").Value = 80005
$ws.Range("M10").Value = -230
$ws.Range("N10").Value = -80345
$ws.Range("H97").Value = 1466.3334
$ws.Range("I97").Value = 1407.8334
$ws.Range("J97").Value = 1583.3334
$ws.Range("K97").Value = 1407.8334
$ws.Range("L97").Value = 1583.3334
$ws.Range("M97").Value = -911.8334
$ws.Range("N97").Value = -2575.3334
$ws.Range("H122").Value = 31300.117
$ws.Range("I122").Value = 44982.78
$ws.Range("J122").Value = 2690.9092
$ws.Range("K122").Value = 134948.34
$ws.Range("L122").Value = 8072.7276
$ws.Range("M122").Value = -132498.34
$ws.Range("N122").Value = -12972.7276
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1702.5
$ws.Range("J99").Value = 2200
$ws.Range("L99").Value = 2200
$ws.Range("N99").Value = -5196
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2306.4
$ws.Range("I16").Value = 1500
$ws.Range("J16").Value = 2508
$ws.Range("K16").Value = 1500
$ws.Range("L16").Value = 2508
$ws.Range("M16").Value = -1213
$ws.Range("N16").Value = -3082
$ws.Range("H31").Value = 6931.553
$ws.Range("I31").Value = 2308.4119
$ws.Range("J31").Value = 9551.333000000001
$ws.Range("K31").Value = 2308.4119
$ws.Range("L31").Value = 9551.333000000001
$ws.Range("M31").Value = -2013.4119
$ws.Range("N31").Value = -10141.333
$ws.Range("H34").Value = 6931.553
$ws.Range("I34").Value = 2308.4119
$ws.Range("J34").Value = 9551.333000000001
$ws.Range("K34").Value = 2308.4119
$ws.Range("L34").Value = 9551.333000000001
$ws.Range("M34").Value = -2106.4119
$ws.Range("N34").Value = -9955.333000000001
$ws.Range("H107").Value = 1036.25
$ws.Range("I107").Value = 685
$ws.Range("J107").Value = 1387.5
$ws.Range("K107").Value = 685
$ws.Range("L107").Value = 1387.5
$ws.Range("M107").Value = 1235
$ws.Range("N107").Value = -5227.5
$ws.Range("H113").Value = 2306.4
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 2508
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 2508
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -6848
$ws.Range("H122").Value = 1751.8966
$ws.Range("I122").Value = 1445
$ws.Range("J122").Value = 1890
$ws.Range("K122").Value = 4335
$ws.Range("L122").Value = 5670
$ws.Range("M122").Value = -1885
$ws.Range("N122").Value = -10570
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 41361.41
$ws.Range("I9").Value = 316.66666
$ws.Range("J9").Value = 47842.156
$ws.Range("K9").Value = 949.9999799999999
$ws.Range("L9").Value = 143526.468
$ws.Range("M9").Value = -725.9999799999999
$ws.Range("N9").Value = -143974.468
$ws.Range("H19").Value = 3627.3333
$ws.Range("J19").Value = 8882
$ws.Range("L19").Value = 26646
$ws.Range("N19").Value = -26994
$ws.Range("H25").Value = 456.25
$ws.Range("J25").Value = 480
$ws.Range("L25").Value = 1440
$ws.Range("N25").Value = -1778
$ws.Range("H30").Value = 456.25
$ws.Range("J30").Value = 480
$ws.Range("L30").Value = 1440
$ws.Range("N30").Value = -1644
$ws.Range("H46").Value = 1509.5834
$ws.Range("I46").Value = 467
$ws.Range("J46").Value = 2254.2856
$ws.Range("K46").Value = 1401
$ws.Range("L46").Value = 6762.8568
$ws.Range("M46").Value = -1310
$ws.Range("N46").Value = -6944.8568
$ws.Range("H70").Value = 6759.5713
$ws.Range("I70").Value = 4663.4
$ws.Range("K70").Value = 13990.2
$ws.Range("M70").Value = -13675.2
$ws.Range("H73").Value = 6759.5713
$ws.Range("I73").Value = 4663.4
$ws.Range("K73").Value = 13990.2
$ws.Range("M73").Value = -12898.2
$ws.Range("H74").Value = 2833.3333
$ws.Range("I74").Value = 1000
$ws.Range("J74").Value = 3750
$ws.Range("K74").Value = 3000
$ws.Range("L74").Value = 11250
$ws.Range("M74").Value = -1939
$ws.Range("N74").Value = -13372
$ws.Range("H77").Value = 2833.3333
$ws.Range("I77").Value = 1000
$ws.Range("J77").Value = 3750
$ws.Range("K77").Value = 9000
$ws.Range("L77").Value = 33750
$ws.Range("M77").Value = -3696
$ws.Range("N77").Value = -44358
$ws.Range("H108").Value = 2451.4
$ws.Range("I108").Value = 509
$ws.Range("K108").Value = 1527
$ws.Range("M108").Value = 1353
$ws.Range("H113").Value = 509.75
$ws.Range("I113").Value = 514.4
$ws.Range("J113").Value = 502
$ws.Range("K113").Value = 1543.2
$ws.Range("L113").Value = 1506
$ws.Range("M113").Value = 626.8000000000002
$ws.Range("N113").Value = -5846
$ws.Range("H114").Value = 3007.4
$ws.Range("I114").Value = 182
$ws.Range("J114").Value = 7245.5
$ws.Range("K114").Value = 546
$ws.Range("L114").Value = 21736.5
$ws.Range("M114").Value = 2708
$ws.Range("N114").Value = -28244.5
$ws.Range("H117").Value = 700
$ws.Range("J117").Value = 700
$ws.Range("L117").Value = 2100
$ws.Range("N117").Value = -8984
$ws.Range("H121").Value = 1904.8182
$ws.Range("I121").Value = 666.6667
$ws.Range("J121").Value = 3390.6
$ws.Range("K121").Value = 2000.0001
$ws.Range("L121").Value = 10171.8
$ws.Range("M121").Value = -690.0001
$ws.Range("N121").Value = -12791.8
$ws.Range("H140").Value = 2043.9375
$ws.Range("I140").Value = 1777
$ws.Range("J140").Value = 2488.8333
$ws.Range("K140").Value = 5331
$ws.Range("L140").Value = 7466.499899999999
$ws.Range("M140").Value = -151
$ws.Range("N140").Value = -17826.4999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2318.182
$ws.Range("I122").Value = 2300
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 6900
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -4450
$ws.Range("N122").Value = -12400
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4580.923
$ws.Range("I7").Value = 4128.9165
$ws.Range("J7").Value = 10005
$ws.Range("K7").Value = 4128.9165
$ws.Range("L7").Value = 10005
$ws.Range("M7").Value = -4016.9165
$ws.Range("N7").Value = -10229
$ws.Range("H41").Value = 10000
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 10000
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 10000
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -10876
$ws.Range("H68").Value = 1634.7142
$ws.Range("I68").Value = 1558.6
$ws.Range("J68").Value = 1825
$ws.Range("K68").Value = 1558.6
$ws.Range("L68").Value = 1825
$ws.Range("M68").Value = -809.5999999999999
$ws.Range("N68").Value = -3323
$ws.Range("H71").Value = 1634.7142
$ws.Range("I71").Value = 1558.6
$ws.Range("J71").Value = 1825
$ws.Range("K71").Value = 7793
$ws.Range("L71").Value = 9125
$ws.Range("M71").Value = -4049
$ws.Range("N71").Value = -16613
$ws.Range("H122").Value = 1473.8
$ws.Range("I122").Value = 1092.25
$ws.Range("K122").Value = 3276.75
$ws.Range("M122").Value = -826.75
$ws.Range("H126").Value = 4580.923
$ws.Range("I126").Value = 4128.9165
$ws.Range("J126").Value = 10005
$ws.Range("K126").Value = 12386.7495
$ws.Range("L126").Value = 30015
$ws.Range("M126").Value = -9916.749500000002
$ws.Range("N126").Value = -34955
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 15000
$ws.Range("J48").Value = 15000
$ws.Range("L48").Value = 15000
$ws.Range("N48").Value = -16138
$ws.Range("H122").Value = 1803.8649
$ws.Range("I122").Value = 1955.92
$ws.Range("J122").Value = 1487.0834
$ws.Range("K122").Value = 5867.76
$ws.Range("L122").Value = 4461.2502
$ws.Range("M122").Value = -3417.76
$ws.Range("N122").Value = -9361.2502
